$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting the existing "posiadane pieniądze" column to E
$ws.Columns("D").Insert()

# New header for the inserted column
$ws.Range("D1").Value = "Koronawirus"

# All rows tested positive
$ws.Range("D2").Value = "pozytywny"
$ws.Range("D3").Value = "pozytywny"
$ws.Range("D4").Value = "pozytywny"
$ws.Range("D5").Value = "pozytywny"

# Updated "posiadane pieniądze" values (now in column E)
$ws.Range("E2").Value = 23
$ws.Range("E3").Value = 220
$ws.Range("E4").Value = 350
$ws.Range("E5").Value = 23
